$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Register new shared strings in the same order the target workbook uses
# (Montpellier, then "Id du trajet", then "null") so the shared-string
# table comes out in the expected order.
$ws.Range("J4").Value = "Montpellier"
$ws.Range("A1").Value = "Id du trajet"
$ws.Range("A3").Value = "null"

# New data cell in row 2
$ws.Range("A2").Value = 12

# New row 4 (full record appended to the "bdd")
$ws.Range("A4").Value = 34
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 44114.416666666664
$ws.Range("E4").Value = 44114.416666666664
$ws.Range("F4").Value = 43545.635775462964
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = "Marseille"
$ws.Range("K4").Value = 1

# Give the D:F date columns on the new row the same date number format as
# the rest of the table (matches styles.xml numFmtId 164 exactly).
$ws.Range("D4:F4").NumberFormat = "yyyy/mm/dd\ hh:mm:ss"

# Column A shrinks to fit the (now shorter) header "Id du trajet"
$ws.Columns.Item(1).AutoFit() | Out-Null

# Move the active selection, like the author's workbook shows
$ws.Range("D9").Select() | Out-Null
